# Swap the content of three pairs of rows (2<->3, 23<->24, 25<->26).
# These rows were re-ordered/re-matched against their GPS points, so every
# populated cell in a pair trades places with its counterpart in the other
# row of the pair (the row number itself - column A's "Id" - moves too).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the columns whose value actually differs between the two rows of at
# least one swapped pair: A, B, E, F, G, H, M, Q, R, Z, AB, AC, AE.
# (Columns that are identical between the paired rows - e.g. the Y/AA date
# columns, which Excel would otherwise "helpfully" reparse as real dates
# when written back through .Value - are intentionally left untouched.)
$cols = @(1,2,5,6,7,8,13,17,18,26,28,29,31)

$pairs = @(@(2,3), @(23,24), @(25,26))

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $vals1 = @{}
    $vals2 = @{}

    foreach ($c in $cols) {
        $vals1[$c] = $ws.Cells.Item($r1, $c).Value()
        $vals2[$c] = $ws.Cells.Item($r2, $c).Value()
    }

    foreach ($c in $cols) {
        $ws.Cells.Item($r1, $c).Value = $vals2[$c]
        $ws.Cells.Item($r2, $c).Value = $vals1[$c]
    }
}
